# Weekly update: insert two new rows (565-566) at the top of the "Feria
# Lagunitas de Puerto Montt - Limon" data block, shifting the existing
# rows 565:603 down to 567:605.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("565:566").Insert()

# New row 565: 1a amarillo
$ws.Cells.Item(565, 1).Value = 4
$ws.Cells.Item(565, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(565, 3).Value = "Los Lagos"
$ws.Cells.Item(565, 4).Value = 44826
$ws.Cells.Item(565, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(565, 5).Value = 10
$ws.Cells.Item(565, 6).Value = "Fruta"
$ws.Cells.Item(565, 7).Value = 100102
$ws.Cells.Item(565, 8).Value = "Cítricos"
$ws.Cells.Item(565, 9).Value = 100102003
$ws.Cells.Item(565, 10).Value = "Limón"
$ws.Cells.Item(565, 11).Value = "Sin especificar"
$ws.Cells.Item(565, 12).Value = "1a amarillo"
$ws.Cells.Item(565, 13).Value = 1000
$ws.Cells.Item(565, 14).Value = 9500
$ws.Cells.Item(565, 15).Value = 10000
$ws.Cells.Item(565, 16).Value = 9750
$ws.Cells.Item(565, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(565, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(565, 19).Value = 542
$ws.Cells.Item(565, 20).Value = 18

# New row 566: 2a amarillo
$ws.Cells.Item(566, 1).Value = 4
$ws.Cells.Item(566, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(566, 3).Value = "Los Lagos"
$ws.Cells.Item(566, 4).Value = 44826
$ws.Cells.Item(566, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(566, 5).Value = 10
$ws.Cells.Item(566, 6).Value = "Fruta"
$ws.Cells.Item(566, 7).Value = 100102
$ws.Cells.Item(566, 8).Value = "Cítricos"
$ws.Cells.Item(566, 9).Value = 100102003
$ws.Cells.Item(566, 10).Value = "Limón"
$ws.Cells.Item(566, 11).Value = "Sin especificar"
$ws.Cells.Item(566, 12).Value = "2a amarillo"
$ws.Cells.Item(566, 13).Value = 400
$ws.Cells.Item(566, 14).Value = 8500
$ws.Cells.Item(566, 15).Value = 8500
$ws.Cells.Item(566, 16).Value = 8500
$ws.Cells.Item(566, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(566, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(566, 19).Value = 472
$ws.Cells.Item(566, 20).Value = 18
